# Add the two new columns (TransactionSpeedNS, TransactionSpeedMS) that were
# computed for the cleaned payments data: the duration between AttemptTime
# (column C) and ResolveTime (column D), in nanoseconds and in rounded
# milliseconds.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Cells.Item(1, 7).Value = "TransactionSpeedNS"
$ws.Cells.Item(1, 8).Value = "TransactionSpeedMS"

# Find the last used row on the sheet (data starts on row 2)
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $attemptTime = $ws.Cells.Item($r, 3).Value2
    $resolveTime = $ws.Cells.Item($r, 4).Value2

    $speedNs = $resolveTime - $attemptTime
    $speedMs = [Math]::Round($speedNs / 1000000.0, 0)

    $ws.Cells.Item($r, 7).Value = $speedNs
    $ws.Cells.Item($r, 8).Value = $speedMs
}
